$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.431604385375977
$ws.Range("B1").Value = 2.406329393386841
$ws.Range("C1").Value = 3.016803741455078
$ws.Range("D1").Value = 3.50059986114502
$ws.Range("E1").Value = 1.884233355522156
